$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("M2").Value = -85.33332999999999
$ws.Range("K2").Value = 198.33333
$ws.Range("I2").Value = 198.33333
$ws.Range("N2").ClearContents()
$ws.Range("H2").Value = 198.33333
$ws.Range("M33").Value = 38.78572
$ws.Range("K33").Value = 190.21428
$ws.Range("I33").Value = 190.21428
$ws.Range("H33").Value = 179
$ws.Range("K52").Value = 3000000
$ws.Range("M52").Value = -2999840
$ws.Range("H52").Value = 38749.04
$ws.Range("I52").Value = 1000000
$ws.Range("K127").Value = 1004.50002
$ws.Range("I127").Value = 334.83334
$ws.Range("H127").Value = 1150.875
$ws.Range("M127").Value = 3955.49998
$ws.Range("L135").Value = 54304.713
$ws.Range("J135").Value = 6033.857
$ws.Range("H135").Value = 2454.3076
$ws.Range("N135").Value = -59374.713
$ws.Range("H137").Value = 4619.5713
$ws.Range("I137").Value = 4812.643
$ws.Range("M137").Value = -11887.929
$ws.Range("K137").Value = 14437.929

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3101.6394
$ws.Range("K32").Value = 3003.8965
$ws.Range("I32").Value = 3003.8965
$ws.Range("M32").Value = -2716.8965
$ws.Range("H37").Value = 57016.375
$ws.Range("J37").Value = 56019.6
$ws.Range("L37").Value = 56019.6
$ws.Range("N37").Value = -56565.6
$ws.Range("J44").Value = 22262.25
$ws.Range("H44").Value = 22262.25
$ws.Range("N44").Value = -23238.25
$ws.Range("L44").Value = 22262.25
$ws.Range("I45").Value = 2990.7
$ws.Range("L45").Value = 13236.375
$ws.Range("H45").Value = 7544.3335
$ws.Range("N45").Value = -13990.375
$ws.Range("J45").Value = 13236.375
$ws.Range("K45").Value = 2990.7
$ws.Range("M45").Value = -2613.7
$ws.Range("J55").Value = 34526.5
$ws.Range("L55").Value = 34526.5
$ws.Range("H55").Value = 36351
$ws.Range("N55").Value = -35156.5
$ws.Range("H74").Value = 22703.23
$ws.Range("K74").Value = 37764.145
$ws.Range("J74").Value = 5132.1665
$ws.Range("I74").Value = 37764.145
$ws.Range("L74").Value = 5132.1665
$ws.Range("M74").Value = -36890.145
$ws.Range("N74").Value = -6880.1665
$ws.Range("J77").Value = 5132.1665
$ws.Range("I77").Value = 37764.145
$ws.Range("L77").Value = 25660.8325
$ws.Range("H77").Value = 22703.23
$ws.Range("K77").Value = 188820.725
$ws.Range("M77").Value = -184452.725
$ws.Range("N77").Value = -34396.8325
$ws.Range("H80").Value = 61565.832
$ws.Range("L80").Value = 72348.75
$ws.Range("N80").Value = -74344.75
$ws.Range("J80").Value = 72348.75
$ws.Range("H83").Value = 61565.832
$ws.Range("J83").Value = 72348.75
$ws.Range("L83").Value = 217046.25
$ws.Range("N83").Value = -227030.25
$ws.Range("M102").Value = -961.1819999999998
$ws.Range("I102").Value = 2583.182
$ws.Range("K102").Value = 2583.182
$ws.Range("H102").Value = 3039.6155
$ws.Range("I132").Value = 3208.9768
$ws.Range("K132").Value = 9626.930399999999
$ws.Range("M132").Value = -7096.930399999999
$ws.Range("L132").Value = 43332.999
$ws.Range("H132").Value = 5153.5576
$ws.Range("N132").Value = -48392.999
$ws.Range("J132").Value = 14444.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2467.8462
$ws.Range("L99").Value = 3179.2
$ws.Range("N99").Value = -6175.2
$ws.Range("M99").Value = -525.25
$ws.Range("K99").Value = 2023.25
$ws.Range("J99").Value = 3179.2
$ws.Range("I99").Value = 2023.25
$ws.Range("I105").Value = 101636.3
$ws.Range("K105").Value = 101636.3
$ws.Range("M105").Value = -99889.3
$ws.Range("H105").Value = 42035.88

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 481683.38
$ws.Range("I58").Value = 2503036.5
$ws.Range("K58").Value = 2503036.5
$ws.Range("M58").Value = -2502833.5
$ws.Range("N86").Value = -39746
$ws.Range("H86").Value = 15184.223
$ws.Range("K86").Value = 4026.3333
$ws.Range("J86").Value = 37500
$ws.Range("M86").Value = -2903.3333
$ws.Range("I86").Value = 4026.3333
$ws.Range("L86").Value = 37500
$ws.Range("L89").Value = 187500
$ws.Range("N89").Value = -198732
$ws.Range("K89").Value = 20131.6665
$ws.Range("H89").Value = 15184.223
$ws.Range("I89").Value = 4026.3333
$ws.Range("J89").Value = 37500
$ws.Range("M89").Value = -14515.6665
$ws.Range("H99").Value = 5908.0835
$ws.Range("L99").Value = 6310.8887
$ws.Range("N99").Value = -9306.8887
$ws.Range("M99").Value = -3201.6665
$ws.Range("K99").Value = 4699.6665
$ws.Range("J99").Value = 6310.8887
$ws.Range("I99").Value = 4699.6665
$ws.Range("H122").Value = 1954.0714
$ws.Range("J122").Value = 3081.8
$ws.Range("L122").Value = 9245.400000000001
$ws.Range("N122").Value = -14145.4
$ws.Range("N126").Value = -23872.6661
$ws.Range("H126").Value = 5908.0835
$ws.Range("J126").Value = 6310.8887
$ws.Range("M126").Value = -11628.9995
$ws.Range("L126").Value = 18932.6661
$ws.Range("I126").Value = 4699.6665
$ws.Range("K126").Value = 14098.9995
$ws.Range("I132").Value = 2765.577
$ws.Range("K132").Value = 8296.731
$ws.Range("M132").Value = -5766.731
$ws.Range("H132").Value = 3590.353
$ws.Range("I136").Value = 2503036.5
$ws.Range("M136").Value = -7506559.5
$ws.Range("H136").Value = 481683.38
$ws.Range("K136").Value = 7509109.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N23").Value = -865.49999
$ws.Range("H23").Value = 158.77777
$ws.Range("J23").Value = 131.83333
$ws.Range("L23").Value = 395.49999
$ws.Range("K108").Value = 42729.999
$ws.Range("M108").Value = -39849.999
$ws.Range("L108").Value = 10500
$ws.Range("H108").Value = 11557.5
$ws.Range("N108").Value = -16260
$ws.Range("I108").Value = 14243.333
$ws.Range("J108").Value = 3500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L11").Value = 15009143
$ws.Range("K11").Value = 6667333.5
$ws.Range("N11").Value = -15009421
$ws.Range("I11").Value = 6667333.5
$ws.Range("H11").Value = 12506600
$ws.Range("J11").Value = 15009143
$ws.Range("M11").Value = -6667194.5
$ws.Range("H18").Value = 12266.667
$ws.Range("N18").Value = -13486
$ws.Range("J18").Value = 12900
$ws.Range("L18").Value = 12900
$ws.Range("J59").Value = 19999
$ws.Range("L59").Value = 19999
$ws.Range("H59").Value = 12499.5
$ws.Range("I59").Value = 5000
$ws.Range("K59").Value = 5000
$ws.Range("M59").Value = -4417
$ws.Range("N59").Value = -21165
$ws.Range("H80").Value = 6004404
$ws.Range("K80").Value = 6670003
$ws.Range("M80").Value = -6669005
$ws.Range("L80").Value = 5006006
$ws.Range("I80").Value = 6670003
$ws.Range("N80").Value = -5008002
$ws.Range("J80").Value = 5006006
$ws.Range("M83").Value = -33345023
$ws.Range("K83").Value = 33350015
$ws.Range("H83").Value = 6004404
$ws.Range("J83").Value = 5006006
$ws.Range("L83").Value = 25030030
$ws.Range("N83").Value = -25040014
$ws.Range("I83").Value = 6670003
$ws.Range("N101").ClearContents()
$ws.Range("J101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("L101").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N23").Value = -53793.332
$ws.Range("K23").Value = 63050
$ws.Range("M23").Value = -62820
$ws.Range("H23").Value = 57220
$ws.Range("J23").Value = 53333.332
$ws.Range("L23").Value = 53333.332
$ws.Range("I23").Value = 63050
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("I25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("L25").Value = 0
$ws.Range("M40").Value = -3334865.2
$ws.Range("H40").Value = 2003900.8
$ws.Range("I40").Value = 3335001.2
$ws.Range("K40").Value = 3335001.2
$ws.Range("H122").Value = 1914727.9
$ws.Range("M122").Value = -15006556
$ws.Range("I122").Value = 5003002
$ws.Range("K122").Value = 15009006
$ws.Range("I132").Value = 2001
$ws.Range("K132").Value = 6003
$ws.Range("M132").Value = -3473
$ws.Range("L132").Value = 24039.375
$ws.Range("H132").Value = 4830.2354
$ws.Range("N132").Value = -29099.375
$ws.Range("J132").Value = 8013.125
$ws.Range("I136").Value = 3911.5
$ws.Range("M136").Value = -9184.5
$ws.Range("H136").Value = 4310.1113
$ws.Range("K136").Value = 11734.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N3").Value = -378
$ws.Range("J3").Value = 150
$ws.Range("H3").Value = 150
$ws.Range("L3").Value = 150
$ws.Range("H54").Value = 20032
$ws.Range("J54").Value = 20032
$ws.Range("L54").Value = 20032
$ws.Range("N54").Value = -21072
$ws.Range("L81").Value = 19898
$ws.Range("J81").Value = 9949
$ws.Range("H81").Value = 6224.5
$ws.Range("N81").Value = -22020
$ws.Range("N84").Value = -110098
$ws.Range("J84").Value = 9949
$ws.Range("L84").Value = 99490
$ws.Range("H84").Value = 6224.5
$ws.Range("I136").Value = 1114038.5
$ws.Range("M136").Value = -3339565.5
$ws.Range("H136").Value = 1003435.2
$ws.Range("K136").Value = 3342115.5
